# Fruta / hortaliza, semanal
# Insert a new weekly record at row 24, pushing the existing rows 24-80 down
# to rows 25-81 (matching the dates/values already present right below them),
# and populate the newly inserted row with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 24; this shifts rows 24:80 -> 25:81
# and carries along the existing row formatting (e.g. the date style on column D).
$ws.Rows.Item(24).Insert()

# Populate the new row 24 with the new data point.
$ws.Range("A24").Value2 = 4
$ws.Range("B24").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C24").Value2 = "Los Lagos"
$ws.Range("D24").Value2 = 44715
$ws.Range("E24").Value2 = 10
$ws.Range("F24").Value2 = 100112031
$ws.Range("G24").Value2 = "Poroto verde"
$ws.Range("H24").Value2 = "Magnum"
$ws.Range("I24").Value2 = "Primera"
$ws.Range("J24").Value2 = 40
$ws.Range("K24").Value2 = 26000
$ws.Range("L24").Value2 = 26000
$ws.Range("M24").Value2 = 26000
$ws.Range("N24").Value2 = "$/malla 25 kilos"
$ws.Range("O24").Value2 = "Perú"
$ws.Range("P24").Value2 = 1040
$ws.Range("Q24").Value2 = 25
$ws.Range("R24").Value2 = "Hortaliza"
